# Auto-generated edit script: apply scheduled-runner data refresh to Jenova_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific leve rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 63762.188
$ws.Range("I28").Value = 63762.188
$ws.Range("K28").Value = 63762.188
$ws.Range("M28").Value = -63277.188

$ws.Range("H41").Value = 35717904
$ws.Range("I41").Value = 998.1667
$ws.Range("J41").Value = 62505580
$ws.Range("K41").Value = 998.1667
$ws.Range("L41").Value = 62505580
$ws.Range("M41").Value = -558.1667
$ws.Range("N41").Value = -62506460

$ws.Range("H76").Value = 83339550
$ws.Range("I76").Value = 5637
$ws.Range("J76").Value = 111117520
$ws.Range("K76").Value = 5637
$ws.Range("L76").Value = 111117520
$ws.Range("M76").Value = -5322
$ws.Range("N76").Value = -111118150

$ws.Range("H79").Value = 83339550
$ws.Range("I79").Value = 5637
$ws.Range("J79").Value = 111117520
$ws.Range("K79").Value = 5637
$ws.Range("L79").Value = 111117520
$ws.Range("M79").Value = -4545
$ws.Range("N79").Value = -111119704

$ws.Range("H138").Value = 3556.9312
$ws.Range("I138").Value = 1493.8
$ws.Range("J138").Value = 4172.791
$ws.Range("K138").Value = 4481.4
$ws.Range("L138").Value = 12518.373
$ws.Range("M138").Value = 658.6000000000004
$ws.Range("N138").Value = -22798.373

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5105.34
$ws.Range("I32").Value = 4808.9165
$ws.Range("K32").Value = 4808.9165
$ws.Range("M32").Value = -4521.9165

$ws.Range("H102").Value = 2242
$ws.Range("I102").Value = 2270.875
$ws.Range("J102").Value = 2011
$ws.Range("K102").Value = 2270.875
$ws.Range("L102").Value = 2011
$ws.Range("M102").Value = -648.875
$ws.Range("N102").Value = -5255

$ws.Range("H122").Value = 5987.3213
$ws.Range("J122").Value = 5478.8125
$ws.Range("L122").Value = 16436.4375
$ws.Range("N122").Value = -21336.4375

$ws.Range("H132").Value = 3143.577
$ws.Range("I132").Value = 3048.353
$ws.Range("K132").Value = 9145.059000000001
$ws.Range("M132").Value = -6615.059000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 27149.684
$ws.Range("I134").Value = 2736.8918
$ws.Range("K134").Value = 8210.6754
$ws.Range("M134").Value = -5675.6754

$ws.Range("H141").Value = 51850
$ws.Range("J141").Value = 44700
$ws.Range("L141").Value = 44700
$ws.Range("N141").Value = -55060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2142
$ws.Range("I16").Value = 1861
$ws.Range("K16").Value = 1861
$ws.Range("M16").Value = -1574

$ws.Range("H86").Value = 8683
$ws.Range("I86").Value = 8045.4
$ws.Range("J86").Value = 10277
$ws.Range("K86").Value = 8045.4
$ws.Range("L86").Value = 10277
$ws.Range("M86").Value = -6922.4
$ws.Range("N86").Value = -12523

$ws.Range("H89").Value = 8683
$ws.Range("I89").Value = 8045.4
$ws.Range("J89").Value = 10277
$ws.Range("K89").Value = 40227
$ws.Range("L89").Value = 51385
$ws.Range("M89").Value = -34611
$ws.Range("N89").Value = -62617

$ws.Range("H105").Value = 1008.875
$ws.Range("I105").Value = 1008.3333
$ws.Range("K105").Value = 1008.3333
$ws.Range("M105").Value = 738.6667

$ws.Range("H113").Value = 2142
$ws.Range("I113").Value = 1861
$ws.Range("K113").Value = 1861
$ws.Range("M113").Value = 309

$ws.Range("H132").Value = 2350
$ws.Range("I132").Value = 2329.5
$ws.Range("K132").Value = 6988.5
$ws.Range("M132").Value = -4458.5

$ws.Range("H138").Value = 49998.168
$ws.Range("J138").Value = 49998.168
$ws.Range("L138").Value = 49998.168
$ws.Range("N138").Value = -60278.168

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1824.75
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 1824.75
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 5474.25
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -6394.25

$ws.Range("H57").Value = 1500
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 1500
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 4500
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -5618

$ws.Range("H82").Value = 7395
$ws.Range("J82").Value = 6960.8335
$ws.Range("L82").Value = 20882.5005
$ws.Range("N82").Value = -21694.5005

$ws.Range("H85").Value = 7395
$ws.Range("J85").Value = 6960.8335
$ws.Range("L85").Value = 20882.5005
$ws.Range("N85").Value = -23690.5005

$ws.Range("H92").Value = 770147.7
$ws.Range("I92").Value = 1428974.6
$ws.Range("J92").Value = 1516.3334
$ws.Range("K92").Value = 4286923.800000001
$ws.Range("L92").Value = 4549.0002
$ws.Range("M92").Value = -4285675.800000001
$ws.Range("N92").Value = -7045.0002

$ws.Range("H131").Value = 12423951
$ws.Range("I131").Value = 30395224
$ws.Range("J131").Value = 68700.81
$ws.Range("K131").Value = 91185672
$ws.Range("L131").Value = 206102.43
$ws.Range("M131").Value = -91180632
$ws.Range("N131").Value = -216182.43

$ws.Range("H132").Value = 462531.16
$ws.Range("I132").Value = 112123.336
$ws.Range("J132").Value = 672775.9
$ws.Range("K132").Value = 1009110.024
$ws.Range("L132").Value = 6054983.100000001
$ws.Range("M132").Value = -1006580.024
$ws.Range("N132").Value = -6060043.100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 19416.666
$ws.Range("J55").Value = 19122.5
$ws.Range("L55").Value = 19122.5
$ws.Range("N55").Value = -19776.5

$ws.Range("H70").Value = 30008.5
$ws.Range("I70").Value = 20008
$ws.Range("J70").Value = 40009
$ws.Range("K70").Value = 20008
$ws.Range("L70").Value = 40009
$ws.Range("M70").Value = -19738
$ws.Range("N70").Value = -40549

$ws.Range("H73").Value = 30008.5
$ws.Range("I73").Value = 20008
$ws.Range("J73").Value = 40009
$ws.Range("K73").Value = 20008
$ws.Range("L73").Value = 40009
$ws.Range("M73").Value = -19072
$ws.Range("N73").Value = -41881

$ws.Range("H102").Value = 2266.5
$ws.Range("I102").Value = 1173.1
$ws.Range("K102").Value = 1173.1
$ws.Range("M102").Value = 448.9000000000001

$ws.Range("H132").Value = 44460.348
$ws.Range("I132").Value = 6287.7827
$ws.Range("J132").Value = 337116.66
$ws.Range("K132").Value = 18863.3481
$ws.Range("L132").Value = 1011349.98
$ws.Range("M132").Value = -16333.3481
$ws.Range("N132").Value = -1016409.98

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10691.833
$ws.Range("I7").Value = 12200.25
$ws.Range("J7").Value = 7675
$ws.Range("K7").Value = 12200.25
$ws.Range("L7").Value = 7675
$ws.Range("M7").Value = -12088.25
$ws.Range("N7").Value = -7899

$ws.Range("H17").Value = 1777.5
$ws.Range("I17").Value = 555
$ws.Range("J17").Value = 3000
$ws.Range("K17").Value = 555
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = -385
$ws.Range("N17").Value = -3340

$ws.Range("H40").Value = 4548.7144
$ws.Range("I40").Value = 3179.2
$ws.Range("J40").Value = 7972.5
$ws.Range("K40").Value = 3179.2
$ws.Range("L40").Value = 7972.5
$ws.Range("M40").Value = -3043.2
$ws.Range("N40").Value = -8244.5

$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992

$ws.Range("H93").Value = 76926000
$ws.Range("J93").Value = 3608.4285
$ws.Range("L93").Value = 3608.4285
$ws.Range("N93").Value = -6104.4285

$ws.Range("H122").Value = 3691
$ws.Range("I122").Value = 3471.6667
$ws.Range("J122").Value = 4349
$ws.Range("K122").Value = 10415.0001
$ws.Range("L122").Value = 13047
$ws.Range("M122").Value = -7965.000100000001
$ws.Range("N122").Value = -17947

$ws.Range("H126").Value = 10691.833
$ws.Range("I126").Value = 12200.25
$ws.Range("J126").Value = 7675
$ws.Range("K126").Value = 36600.75
$ws.Range("L126").Value = 23025
$ws.Range("M126").Value = -34130.75
$ws.Range("N126").Value = -27965

$ws.Range("H132").Value = 6429.7393
$ws.Range("I132").Value = 6194.4
$ws.Range("K132").Value = 18583.2
$ws.Range("M132").Value = -16053.2

$ws.Range("H133").Value = 59766.445
$ws.Range("J133").Value = 59766.445
$ws.Range("L133").Value = 59766.445
$ws.Range("N133").Value = -64826.445

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 3875
$ws.Range("I17").Value = 3875
$ws.Range("K17").Value = 3875
$ws.Range("M17").Value = -3703

$ws.Range("H88").Value = 79592
$ws.Range("J88").Value = 79592
$ws.Range("L88").Value = 79592
$ws.Range("N88").Value = -80404

$ws.Range("H91").Value = 79592
$ws.Range("J91").Value = 79592
$ws.Range("L91").Value = 79592
$ws.Range("N91").Value = -82400

$ws.Range("H107").Value = 817.3913
$ws.Range("I107").Value = 897.55554
$ws.Range("K107").Value = 2692.66662
$ws.Range("M107").Value = -772.66662

$ws.Range("H122").Value = 31251630
$ws.Range("I122").Value = 37038770
$ws.Range("J122").Value = 1089
$ws.Range("K122").Value = 111116310
$ws.Range("L122").Value = 3267
$ws.Range("M122").Value = -111113860
$ws.Range("N122").Value = -8167

$ws.Range("H124").Value = 92463.336
$ws.Range("J124").Value = 92463.336
$ws.Range("L124").Value = 92463.336
$ws.Range("N124").Value = -102283.336

$ws.Range("H132").Value = 45929.543
$ws.Range("I132").Value = 2996.1177
$ws.Range("J132").Value = 150196.42
$ws.Range("K132").Value = 8988.3531
$ws.Range("L132").Value = 450589.26
$ws.Range("M132").Value = -6458.3531
$ws.Range("N132").Value = -455649.26

$ws.Range("H136").Value = 12422737
$ws.Range("I136").Value = 13891954
$ws.Range("J136").Value = 668999
$ws.Range("K136").Value = 41675862
$ws.Range("L136").Value = 2006997
$ws.Range("M136").Value = -41673312
$ws.Range("N136").Value = -2012097
